{"js": "// Sprint Retrospective rewrite: replace the three body paragraphs\n// (the \"went well\", \"issue/tools\" and \"improvement\" paragraphs) with\n// the author's expanded text, and drop the stray \"_GoBack\" bookmark\n// that used to sit at the end of the last paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst wentWellText =\n  \"The thing that went well this spring was developing a clear and concise \" +\n  \"web page that was easy to navigate. As this was my first time working \" +\n  \"with Flask, I was proud of how much progress I made and how the \" +\n  \"development of the website went smoothly. The utilization of ChatGPT \" +\n  \"also helped flesh out the website using AI imaging.\";\n\nconst toolsText =\n  \"Although there were some minor issues during this sprint, there were \" +\n  \"no severe issues that hindered the project\\u2019s progress. Tools such as \" +\n  \"Github\\u2019s project board helped organize the tasks that needed to be \" +\n  \"completed and in what order they should be focused on. Another tool \" +\n  \"that was utilized during this project that mitigates issues was \" +\n  \"ChatGPT. Due to my partner leaving the group without warning, I was \" +\n  \"forced to work and learn Flask on my own. For any tasks that seemed \" +\n  \"complex, I utilized ChatGPT to point me in the right direction in \" +\n  \"solving a specific issue. \";\n\nconst improvementText =\n  \"I do not believe there could be any improvements that could be made \" +\n  \"for this sprint. Although it was a rocky start at the beginning of the \" +\n  \"sprint; being introduced to multiple topics in such a small \" +\n  \"timeframe, I felt we, as a class, were given enough freedom with our \" +\n  \"projects that it was not so overwhelming later down the project\\u2019s \" +\n  \"timeline.\";\n\n// Paragraph 0 is the \"Sprint Retrospective\" title, so the three body\n// paragraphs that get rewritten are items 1, 2 and 3.\nparagraphs.items[1].insertText(wentWellText, \"Replace\");\nparagraphs.items[2].insertText(toolsText, \"Replace\");\nparagraphs.items[3].insertText(improvementText, \"Replace\");\n\n// The old \"An improvement...\" paragraph ended with an empty _GoBack\n// bookmark (left over from the last edit position); the rewritten\n// paragraph no longer has it.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Sprint Retrospective rewrite: replace the three body paragraphs\n# (the \"went well\", \"issue/tools\" and \"improvement\" paragraphs) with\n# the author's expanded text, and drop the stray \"_GoBack\" bookmark\n# that used to sit at the end of the last paragraph.\n\n$d = $word.ActiveDocument\n\n$wentWellText = \"The thing that went well this spring was developing a clear and concise web page that was easy to navigate. As this was my first time working with Flask, I was proud of how much progress I made and how the development of the website went smoothly. The utilization of ChatGPT also helped flesh out the website using AI imaging.\"\n\n$toolsText = \"Although there were some minor issues during this sprint, there were no severe issues that hindered the project\" + [char]0x2019 + \"s progress. Tools such as Github\" + [char]0x2019 + \"s project board helped organize the tasks that needed to be completed and in what order they should be focused on. Another tool that was utilized during this project that mitigates issues was ChatGPT. Due to my partner leaving the group without warning, I was forced to work and learn Flask on my own. For any tasks that seemed complex, I utilized ChatGPT to point me in the right direction in solving a specific issue. \"\n\n$improvementText = \"I do not believe there could be any improvements that could be made for this sprint. Although it was a rocky start at the beginning of the sprint; being introduced to multiple topics in such a small timeframe, I felt we, as a class, were given enough freedom with our projects that it was not so overwhelming later down the project\" + [char]0x2019 + \"s timeline.\"\n\n# Paragraph 1 is the \"Sprint Retrospective\" title, so the three body\n# paragraphs that get rewritten are paragraphs 2, 3 and 4 (1-based COM\n# index). Re-fetch each paragraph's Range right before the assignment\n# and drive the replace through $d.Range(start,end) \u2014 setting\n# $paragraph.Range.Text directly only clobbered the paragraph's first\n# run when the paragraph held several runs.\n$r2 = $d.Paragraphs(2).Range\n$d.Range($r2.Start, $r2.End).Text = $wentWellText\n\n$r3 = $d.Paragraphs(3).Range\n$d.Range($r3.Start, $r3.End).Text = $toolsText\n\n$r4 = $d.Paragraphs(4).Range\n$d.Range($r4.Start, $r4.End).Text = $improvementText\n\n# The old \"An improvement...\" paragraph ended with an empty _GoBack\n# bookmark (left over from the last edit position); the rewritten\n# paragraph no longer has it.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
